# Update cache hit-rate derived values on the "Ways" sheet.
# Columns: C (reqs?), G, H, K, L change per-row; D and I (hit rate %) stay the same.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ways")
$ws.Activate()

$rows = @(
    @{ r = 4;  C = 95191124;  G = 648862226; H = 648859317; K = 11250430; L = 5613899 },
    @{ r = 8;  C = 26625618;  G = 551601635; H = 551600791; K = 11250430; L = 5597781 },
    @{ r = 12; C = 271058;    G = 2812671;   H = 2812578;   K = 11250430; L = 395500 },
    @{ r = 16; C = 81495232;  G = 590940291; H = 590934748; K = 11250430; L = 5613744 },
    @{ r = 20; C = 24509011;  G = 262586425; H = 262584698; K = 11250430; L = 5597453 },
    @{ r = 24; C = 251135;    G = 703295;    H = 703263;    K = 11250430; L = 392073 },
    @{ r = 28; C = 57195124;  G = 624475211; H = 624473369; K = 11250430; L = 5584490 },
    @{ r = 32; C = 15395621;  G = 541889711; H = 541889183; K = 11250430; L = 5331209 },
    @{ r = 36; C = 415054;    G = 3769695;   H = 3769601;   K = 11250430; L = 833168 },
    @{ r = 40; C = 43061838;  G = 435196489; H = 435193716; K = 11250430; L = 5496181 },
    @{ r = 44; C = 5414183;   G = 194486570; H = 194485958; K = 11250430; L = 4627080 },
    @{ r = 48; C = 283048;    G = 926566;    H = 926534;    K = 11250430; L = 513581 }
)

foreach ($row in $rows) {
    $ws.Range("C$($row.r)").Value = $row.C
    $ws.Range("G$($row.r)").Value = $row.G
    $ws.Range("H$($row.r)").Value = $row.H
    $ws.Range("K$($row.r)").Value = $row.K
    $ws.Range("L$($row.r)").Value = $row.L
}

$ws.Range("H13").Select()
